# Trade #7 closed at 2026-02-16 21:21:22 - leadlag UP +0.000%
# Append a new trade row (row 6) to the "leadlag" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 6

$ws.Cells.Item($row, 1).Value = 7                 # A: Trade #
# B: Date - force literal text so "2026-02-16" is not auto-converted to a date serial.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).ClearFormats()
$ws.Cells.Item($row, 3).Value = "21:21:22"         # C: Time
$ws.Cells.Item($row, 4).Value = "leadlag"          # D: Strategy
$ws.Cells.Item($row, 5).Value = "UP"               # E: Side
$ws.Cells.Item($row, 6).Value = 69441.86           # F: Entry Price
# G: Exit Price - still open, leave blank (no exit price yet).
$ws.Cells.Item($row, 8).Value = "OPEN"             # H: Status
$ws.Cells.Item($row, 9).Value = 0                  # I: P&L %
$ws.Cells.Item($row, 10).Value = 0                 # J: P&L $
$ws.Cells.Item($row, 11).Value = 0.75              # K: Confidence
$ws.Cells.Item($row, 12).Value = "Binance leading with 0.079% move"  # L: Entry Reason
# M: Exit Reason - still open, leave blank.
$ws.Cells.Item($row, 14).Value = 0                 # N: Duration (min)
